# Fixed for new file format
# Rename ESC channel labels in column K (yChannel) from "ESC<n>/..." to
# "ESC_<n>/..." to match the new log file naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K6").Value  = "ESC_1/RPM"
$ws.Range("K7").Value  = "ESC_2/RPM"
$ws.Range("K8").Value  = "ESC_3/RPM"
$ws.Range("K9").Value  = "ESC_4/RPM"

$ws.Range("K10").Value = "ESC_1/Temp"
$ws.Range("K11").Value = "ESC_2/Temp"
$ws.Range("K12").Value = "ESC_3/Temp"
$ws.Range("K13").Value = "ESC_4/Temp"

$ws.Range("K14").Value = "ESC_1/Volt"
$ws.Range("K15").Value = "ESC_2/Volt"
$ws.Range("K16").Value = "ESC_3/Volt"
$ws.Range("K17").Value = "ESC_4/Volt"

$ws.Range("K19").Value = "ESC_1/Curr"
$ws.Range("K20").Value = "ESC_2/Curr"
$ws.Range("K21").Value = "ESC_3/Curr"
$ws.Range("K22").Value = "ESC_4/Curr"

# Update the sheet's active cell selection to match the saved view state.
$ws.Range("K23").Select() | Out-Null
